$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03785925370218
$ws.Range("D2").Value = 1.04524604160471
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.053790055485428
$ws.Range("I2").Value = 1.037280911380811
$ws.Range("J2").Value = 1.042959835490737
$ws.Range("K2").Value = 1.048014820969711
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.056535052745712
$ws.Range("N2").Value = 1.018178382498527
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039045880737935
$ws.Range("D3").Value = 1.046164771942056
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.054842783620176
$ws.Range("I3").Value = 1.037519192159759
$ws.Range("J3").Value = 1.043789825735101
$ws.Range("K3").Value = 1.048744673321982
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.057400282218778
$ws.Range("N3").Value = 1.018459634241917
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039813424913413
$ws.Range("D4").Value = 1.046758848677108
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.055523861995926
$ws.Range("I4").Value = 1.03767188201818
$ws.Range("J4").Value = 1.044326090291865
$ws.Range("K4").Value = 1.049215906405099
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.057959435242298
$ws.Range("N4").Value = 1.018641203521407
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040136034459991
$ws.Range("D5").Value = 1.047008502464379
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.055810162324627
$ws.Range("I5").Value = 1.037735715376262
$ws.Range("J5").Value = 1.044551346505638
$ws.Range("K5").Value = 1.049413766704775
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.058194334194198
$ws.Range("N5").Value = 1.018717435032862
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040190198199891
$ws.Range("D6").Value = 1.0470504148134
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.055858231960476
$ws.Range("I6").Value = 1.037746412320644
$ws.Range("J6").Value = 1.044589156932336
$ws.Range("K6").Value = 1.049446973904652
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.05823376486988
$ws.Range("N6").Value = 1.018730228768067
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039817735895279
$ws.Range("D7").Value = 1.046762184940164
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.055527687652471
$ws.Range("I7").Value = 1.037672736367023
$ws.Range("J7").Value = 1.044329100919623
$ws.Range("K7").Value = 1.0492185511904
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.057962574636835
$ws.Range("N7").Value = 1.018642222524125
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038260339453638
$ws.Range("D8").Value = 1.045556614593963
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.054145852445229
$ws.Range("I8").Value = 1.037361748862569
$ws.Range("J8").Value = 1.043240499614611
$ws.Range("K8").Value = 1.048261691625087
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.056827608112215
$ws.Range("N8").Value = 1.018273519799433
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035513765333836
$ws.Range("D9").Value = 1.043429135992441
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.051710033160227
$ws.Range("I9").Value = 1.036802302045913
$ws.Range("J9").Value = 1.041316121389348
$ws.Range("K9").Value = 1.046567667916862
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.054822200287946
$ws.Range("N9").Value = 1.017620597453043
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033681099025833
$ws.Range("D10").Value = 1.042008693196644
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.050085539449362
$ws.Range("I10").Value = 1.036421630169053
$ws.Range("J10").Value = 1.040029030136272
$ws.Range("K10").Value = 1.045432959897615
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.053481555565409
$ws.Range("N10").Value = 1.017183136749386
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032887125701578
$ws.Range("D11").Value = 1.041393113755787
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.049381959425644
$ws.Range("I11").Value = 1.036254964084407
$ws.Range("J11").Value = 1.039470702682946
$ws.Range("K11").Value = 1.044940338614154
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.052900152590287
$ws.Range("N11").Value = 1.016993191236495
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032592143786833
$ws.Range("D12").Value = 1.041164381440408
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.049120592878882
$ws.Range("I12").Value = 1.036192781230252
$ws.Range("J12").Value = 1.039263162271762
$ws.Range("K12").Value = 1.044757163155202
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.052684058212966
$ws.Range("N12").Value = 1.016922558240891
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032655421380816
$ws.Range("D13").Value = 1.041213448873329
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.049176658065433
$ws.Range("I13").Value = 1.036206132142346
$ws.Range("J13").Value = 1.039307687325851
$ws.Range("K13").Value = 1.044796463705321
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.052730417334334
$ws.Range("N13").Value = 1.016937712828254
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032862743740935
$ws.Range("D14").Value = 1.041374208281566
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.049360355313558
$ws.Range("I14").Value = 1.036249829655278
$ws.Range("J14").Value = 1.039453550450695
$ws.Range("K14").Value = 1.044925201235884
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.052882292929902
$ws.Range("N14").Value = 1.016987354299939
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032990473225741
$ws.Range("D15").Value = 1.041473247004988
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.04947353382782
$ws.Range("I15").Value = 1.036276716599517
$ws.Range("J15").Value = 1.039543401260244
$ws.Range("K15").Value = 1.045004494946341
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.052975850522638
$ws.Range("N15").Value = 1.017017929600709
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033733783447653
$ws.Range("D16").Value = 1.042049536185492
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.050132230205188
$ws.Range("I16").Value = 1.036432652609557
$ws.Range("J16").Value = 1.040066063139062
$ws.Range("K16").Value = 1.045465626395776
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.053520122427798
$ws.Range("N16").Value = 1.017195731778571
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034199929147875
$ws.Range("D17").Value = 1.042410887806651
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.050545368533119
$ws.Range("I17").Value = 1.036529976274767
$ws.Range("J17").Value = 1.040393644271744
$ws.Range("K17").Value = 1.045754537253323
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.053861289525906
$ws.Range("N17").Value = 1.017307122363627
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034471783995301
$ws.Range("D18").Value = 1.042621608139253
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.050786329423413
$ws.Range("I18").Value = 1.036586566719821
$ws.Range("J18").Value = 1.040584619587953
$ws.Range("K18").Value = 1.045922930124603
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.054060200237985
$ws.Range("N18").Value = 1.017372044251006
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034564472747121
$ws.Range("D19").Value = 1.042693449878499
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.050868488226467
$ws.Range("I19").Value = 1.036605832619711
$ws.Range("J19").Value = 1.040649720782612
$ws.Range("K19").Value = 1.045980326745442
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.054128009043915
$ws.Range("N19").Value = 1.017394172405321
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034149920286194
$ws.Range("D20").Value = 1.042372123371354
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.050501044333445
$ws.Range("I20").Value = 1.036519552656664
$ws.Range("J20").Value = 1.040358507981656
$ws.Range("K20").Value = 1.045723552679828
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.053824694457664
$ws.Range("N20").Value = 1.017295176419295
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0328016942822
$ws.Range("D21").Value = 1.041326870821919
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.049306261767491
$ws.Range("I21").Value = 1.036236969438663
$ws.Range("J21").Value = 1.039410601622106
$ws.Range("K21").Value = 1.044887296603814
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.052837573098609
$ws.Range("N21").Value = 1.016972738301033
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031953634637219
$ws.Range("D22").Value = 1.040669223295299
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.048554905923779
$ws.Range("I22").Value = 1.036057703173592
$ws.Range("J22").Value = 1.038813730979094
$ws.Range("K22").Value = 1.044360386462797
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.052216146708437
$ws.Range("N22").Value = 1.016769552727087
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032403243541065
$ws.Range("D23").Value = 1.041017898171868
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.048953228306723
$ws.Range("I23").Value = 1.036152886903852
$ws.Range("J23").Value = 1.039130227632474
$ws.Range("K23").Value = 1.044639818221103
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.052545651235811
$ws.Range("N23").Value = 1.016877308571868
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034172517261981
$ws.Range("D24").Value = 1.042389639505566
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.050521072580047
$ws.Range("I24").Value = 1.036524263187351
$ws.Range("J24").Value = 1.040374384859188
$ws.Range("K24").Value = 1.045737553657943
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.053841230460511
$ws.Range("N24").Value = 1.017300574432967
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036224097553
$ws.Range("D25").Value = 1.043979512082622
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.05233985588705
$ws.Range("I25").Value = 1.036948290117123
$ws.Range("J25").Value = 1.041814350637852
$ws.Range("K25").Value = 1.047006555430932
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.055341296512196
$ws.Range("N25").Value = 1.017789776618918
